$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("数学")

# Copy the formatting of the last existing data row (row 4) down onto the
# new row 5 so the new row shares the same cell style (wrap text, etc.)
$ws.Range("A4:G4").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new "sum of squares" problem row.
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 633

$ws.Cells.Item(5, 4).Value = "1 存在两个整数a,b，例如-2，-1，-，1，2，3，4，5，他们的平方之和等于指定的数c`n2 a*a,b*b的值分别小于c，也就是a,b的区间范围是[0,sqrt(c)]`n3 双指针i,j分别指向区间两端，如他们的平方和偏大，j--；如平方和偏小，i++。直到找到这两个数字`n4 循环终止条件j>i"

$ws.Cells.Item(5, 3).Value = "给定一个非负整数 c ，你要判断是否存在两个整数 a 和 b，使得 a2 + b2 = c。 `n输入: 5`n输出: True`n解释: 1 * 1 + 2 * 2 = 5"

$ws.Cells.Item(5, 5).Value = "双指针`n逼近"

$ws.Cells.Item(5, 6).Value = "O(sqrt(n))，n是数字值"

$ws.Cells.Item(5, 7).Value = "O(1)"

# Match the authored row height for the newly added row.
$ws.Rows.Item(5).RowHeight = 154

# Update the selected / active cell shown on this sheet.
$null = $ws.Range("D8").Select()
